# Add test data and reST
# - Rename the existing sheet to "Sheet1"
# - Add a new "Sheet2" worksheet with a small item/qty table
# - Style + merge the item-name column, matching the authored workbook

$wb = $excel.ActiveWorkbook

$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet1"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "H1"
$ws2.Range("B1").Value = "H2"

# Item 1 block (rows 2-4)
$ws2.Range("A2").Value = "Item 1"
$ws2.Range("B2").Value = 1
$ws2.Range("B3").Value = 2
$ws2.Range("B4").Value = 3

# Item 2 block (rows 5-7)
$ws2.Range("A5").Value = "Item 2"
$ws2.Range("B5").Value = 1
$ws2.Range("B6").Value = 2

# Vertically center the item-name column, then merge each item's rows
$ws2.Range("A2:A7").VerticalAlignment = -4108
$ws2.Range("A5:A7").Merge()
$ws2.Range("A2:A4").Merge()

# Vertically center the trailing qty cells of the second item
$ws2.Range("B6:B7").VerticalAlignment = -4108

# Restore the view state: Sheet2 keeps a stale selection from editing,
# but Sheet1 stays the active/selected sheet.
[void]$ws2.Range("C17").Select()
[void]$ws1.Select()
[void]$ws1.Range("C4").Select()
